# This edit re-orders the data rows (rows 2..33) of the sheet according to
# a fixed permutation (the header row 1 and the last data row 34 stay put).
# The used range of the sheet is A..AY.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the full block of data rows (2..33) into a 2-D array (1-based,
# dim1 = row offset 1..32 <-> sheet rows 2..33, dim2 = col offset 1..51 <-> A..AY)
$srcRange = $ws.Range("A2:AY33")
$srcVals = $srcRange.Value()

# Build a new array holding the same data but in the new row order.
# newRowIndex (1..32, sheet row = newRowIndex+1) takes its data from
# srcVals row oldRowIndex (1..32, sheet row = oldRowIndex+1).
$order = @(9,10,11,12,13,14,15,16,1,2,3,17,18,19,4,20,5,21,22,23,24,6,25,26,27,28,29,30,7,31,32,8)

$numCols = $srcVals.GetLength(1)
$newVals = New-Object 'object[,]' 32,$numCols

for ($i = 1; $i -le 32; $i++) {
    $oldRow = $order[$i - 1]
    for ($j = 1; $j -le $numCols; $j++) {
        $newVals[$i - 1, $j - 1] = $srcVals[$oldRow, $j]
    }
}

# Columns A, B, E, Q, R, S hold plain numbers; AD, AE, AG hold booleans.
# Every other column holds text (some of which look like dates or numbers,
# e.g. "2019-11-05" or "1") that Excel would otherwise auto-convert when
# assigned through .Value. Force those columns to text first, then clear
# the format stamp afterwards so the cells end up as plain, unstyled text
# cells (matching their original un-styled state).
# NOTE: NumberFormat only applies to the FIRST area of a multi-area range,
# so each contiguous block is addressed (and formatted) separately.
$textRanges = @(
    $ws.Range("C2:D33"),
    $ws.Range("F2:P33"),
    $ws.Range("T2:AC33"),
    $ws.Range("AF2:AF33"),
    $ws.Range("AH2:AY33")
)
foreach ($r in $textRanges) { $r.NumberFormat = "@" }

$dstRange = $ws.Range("A2:AY33")
$dstRange.Value = $newVals

foreach ($r in $textRanges) { $r.ClearFormats() }
